# Update results table with new metrics for random_forest, lsboost, neural_network
# (recomputed against the updated 2018/2019 dataset) and append a new
# "old_model" row comparing against the old model configuration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same text, but widen column A slightly to fit "old_model"/"random_forest".
# (Target raw width is 15.28515625; the COM ColumnWidth setter here quantizes to
# 1/6-character steps, so 14.45 is the closest input that lands on that bucket.)
$ws.Columns.Item(1).ColumnWidth = 14.45

# Row 2: random_forest
$ws.Range("B2").Value = 4.037626131366113
$ws.Range("C2").Value = 0.2978228181912515
$ws.Range("D2").Value = 3.2223005063056425
$ws.Range("E2").Value = 0.31185920379413629
$ws.Range("F2").Value = 0.55844355470731
$ws.Range("G2").Value = 0.52609051319817945
$ws.Range("H2").Value = 0.68814079620586366
$ws.Range("I2").Value = 0.84187762369780006

# Row 3: lsboost
$ws.Range("B3").Value = 4.3135270726604054
$ws.Range("C3").Value = 0.31817378512193251
$ws.Range("D3").Value = 3.3451288442955018
$ws.Range("E3").Value = 0.35593559144146103
$ws.Range("F3").Value = 0.59660337867083946
$ws.Range("G3").Value = 0.54614414359109698
$ws.Range("H3").Value = 0.64406440855853897
$ws.Range("I3").Value = 0.80786161282634184

# Row 4: neural_network
$ws.Range("B4").Value = 3.9220981055519366
$ws.Range("C4").Value = 0.28930125598895612
$ws.Range("D4").Value = 3.1808095866540169
$ws.Range("E4").Value = 0.29426815493798358
$ws.Range("F4").Value = 0.54246488820750749
$ws.Range("G4").Value = 0.51931647732850394
$ws.Range("H4").Value = 0.70573184506201647
$ws.Range("I4").Value = 0.84473996285894271

# Row 5: old_model (new row)
$ws.Range("A5").Value = "old_model"
$ws.Range("B5").Value = 4.0021211109877983
$ws.Range("C5").Value = 0.29520390180697753
$ws.Range("D5").Value = 3.1078516135057463
$ws.Range("E5").Value = 0.30639862695812736
$ws.Range("F5").Value = 0.55353285987204714
$ws.Range("G5").Value = 0.50740495713963174
$ws.Range("H5").Value = 0.69360137304187264
$ws.Range("I5").Value = 0.84623167537412658
